# MaterialFilter main body: fix composition columns (A-G) that were
# accidentally entered on a 0-100 percentage scale instead of the
# 0-1 fraction scale used throughout the rest of the sheet. Divide
# each affected cell by 100 so it matches the fractional convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "18"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "61"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E")) {
    $addr = $col + "95"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "118"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "166"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "167"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "183"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "203"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "250"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "265"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "269"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "298"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "392"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "408"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "469"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "487"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "488"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "503"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "538"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E")) {
    $addr = $col + "571"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "573"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "583"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "600"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "635"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "703"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "716"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "741"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "779"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "799"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "843"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F")) {
    $addr = $col + "888"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E")) {
    $addr = $col + "902"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E")) {
    $addr = $col + "971"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E")) {
    $addr = $col + "973"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "1011"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
foreach ($col in @("A","B","C","D","E","F","G")) {
    $addr = $col + "1016"
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 / 100
}
